# Apply scheduled-runner price/profit updates to the Belias sheets.
$wb = $excel.ActiveWorkbook

# ALC row 124
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(124, 8).Value = 22666.666
$ws.Cells.Item(124, 10).Value = 22666.666
$ws.Cells.Item(124, 12).Value = 22666.666
$ws.Cells.Item(124, 14).Value = -32486.666

# ALC row 126
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(126, 8).Value = 22400
$ws.Cells.Item(126, 10).Value = 22400
$ws.Cells.Item(126, 12).Value = 22400
$ws.Cells.Item(126, 14).Value = -32280

# ALC row 130
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(130, 8).Value = 25555.555
$ws.Cells.Item(130, 10).Value = 25555.555
$ws.Cells.Item(130, 12).Value = 25555.555
$ws.Cells.Item(130, 14).Value = -35595.555

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 2328.3333
$ws.Cells.Item(137, 9).Value = 2241.5833
$ws.Cells.Item(137, 10).Value = 2501.8333
$ws.Cells.Item(137, 11).Value = 6724.749899999999
$ws.Cells.Item(137, 12).Value = 7505.499899999999
$ws.Cells.Item(137, 13).Value = -4174.749899999999
$ws.Cells.Item(137, 14).Value = -12605.4999

# ARM row 5
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 388.75
$ws.Cells.Item(5, 9).Value = 301.42856
$ws.Cells.Item(5, 10).Value = 1000
$ws.Cells.Item(5, 11).Value = 301.42856
$ws.Cells.Item(5, 12).Value = 1000
$ws.Cells.Item(5, 13).Value = -189.42856
$ws.Cells.Item(5, 14).Value = -1224

# ARM row 123
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(123, 8).Value = 20000
$ws.Cells.Item(123, 10).Value = 20000
$ws.Cells.Item(123, 12).Value = 20000
$ws.Cells.Item(123, 14).Value = -29800

# ARM row 128
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(128, 8).Value = 24638.4
$ws.Cells.Item(128, 10).Value = 24638.4
$ws.Cells.Item(128, 12).Value = 24638.4
$ws.Cells.Item(128, 14).Value = -34598.4

# ARM row 129
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(129, 8).Value = 38749.75
$ws.Cells.Item(129, 10).Value = 38749.75
$ws.Cells.Item(129, 12).Value = 38749.75
$ws.Cells.Item(129, 14).Value = -48749.75

# ARM row 131
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(131, 8).Value = 22990
$ws.Cells.Item(131, 10).Value = 22990
$ws.Cells.Item(131, 12).Value = 22990
$ws.Cells.Item(131, 14).Value = -33070

# BSM row 4
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 388.75
$ws.Cells.Item(4, 9).Value = 301.42856
$ws.Cells.Item(4, 10).Value = 1000
$ws.Cells.Item(4, 11).Value = 301.42856
$ws.Cells.Item(4, 12).Value = 1000
$ws.Cells.Item(4, 13).Value = -186.42856
$ws.Cells.Item(4, 14).Value = -1230

# CRP row 20
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(20, 8).Value = 27499.75
$ws.Cells.Item(20, 10).Value = 27499.75
$ws.Cells.Item(20, 12).Value = 27499.75
$ws.Cells.Item(20, 14).Value = -27971.75

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 283.66666
$ws.Cells.Item(22, 9).Value = 156.66667
$ws.Cells.Item(22, 11).Value = 156.66667
$ws.Cells.Item(22, 13).Value = 193.33333

# CRP row 30
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(30, 8).Value = 27499.75
$ws.Cells.Item(30, 10).Value = 27499.75
$ws.Cells.Item(30, 12).Value = 27499.75
$ws.Cells.Item(30, 14).Value = -27681.75

# CRP row 123
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(123, 8).Value = 21666.666
$ws.Cells.Item(123, 10).Value = 21666.666
$ws.Cells.Item(123, 12).Value = 21666.666
$ws.Cells.Item(123, 14).Value = -31466.666

# CRP row 127
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(127, 8).Value = 31897.334
$ws.Cells.Item(127, 10).Value = 31897.334
$ws.Cells.Item(127, 12).Value = 31897.334
$ws.Cells.Item(127, 14).Value = -41817.334

# CRP row 128
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(128, 8).Value = 27499.75
$ws.Cells.Item(128, 10).Value = 27499.75
$ws.Cells.Item(128, 12).Value = 27499.75
$ws.Cells.Item(128, 14).Value = -37459.75

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 1086.3043
$ws.Cells.Item(122, 9).Value = 491.92307
$ws.Cells.Item(122, 10).Value = 1859
$ws.Cells.Item(122, 11).Value = 4427.30763
$ws.Cells.Item(122, 12).Value = 16731
$ws.Cells.Item(122, 13).Value = -1977.30763
$ws.Cells.Item(122, 14).Value = -21631

# CUL row 123
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(123, 8).Value = 1530
$ws.Cells.Item(123, 9).Value = 1530
$ws.Cells.Item(123, 10).Value = 0
$ws.Cells.Item(123, 11).Value = 4590
$ws.Cells.Item(123, 12).Value = 0
$ws.Cells.Item(123, 13).Value = -2140
$ws.Cells.Item(123, 14).ClearContents()

# CUL row 125
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(125, 8).Value = 7400
$ws.Cells.Item(125, 9).Value = 0
$ws.Cells.Item(125, 10).Value = 7400
$ws.Cells.Item(125, 11).Value = 0
$ws.Cells.Item(125, 12).Value = 22200
$ws.Cells.Item(125, 14).Value = -32040
$ws.Cells.Item(125, 13).ClearContents()

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 865.95
$ws.Cells.Item(131, 10).Value = 910.98865
$ws.Cells.Item(131, 12).Value = 2732.96595
$ws.Cells.Item(131, 14).Value = -12812.96595

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 31978.848
$ws.Cells.Item(122, 9).Value = 38449.406
$ws.Cells.Item(122, 10).Value = 2861.3333
$ws.Cells.Item(122, 11).Value = 115348.218
$ws.Cells.Item(122, 12).Value = 8583.999899999999
$ws.Cells.Item(122, 13).Value = -112898.218
$ws.Cells.Item(122, 14).Value = -13483.9999

# GSM row 124
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(124, 8).Value = 24235.295
$ws.Cells.Item(124, 10).Value = 24235.295
$ws.Cells.Item(124, 12).Value = 24235.295
$ws.Cells.Item(124, 14).Value = -34055.295

# GSM row 128
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(128, 8).Value = 20000
$ws.Cells.Item(128, 10).Value = 20000
$ws.Cells.Item(128, 12).Value = 20000
$ws.Cells.Item(128, 14).Value = -29960

# GSM row 130
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(130, 8).Value = 20000
$ws.Cells.Item(130, 10).Value = 20000
$ws.Cells.Item(130, 12).Value = 20000
$ws.Cells.Item(130, 14).Value = -30040

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 719.087
$ws.Cells.Item(22, 9).Value = 417.5
$ws.Cells.Item(22, 10).Value = 782.5789
$ws.Cells.Item(22, 11).Value = 417.5
$ws.Cells.Item(22, 12).Value = 782.5789
$ws.Cells.Item(22, 13).Value = -122.5
$ws.Cells.Item(22, 14).Value = -1372.5789

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 719.087
$ws.Cells.Item(27, 9).Value = 417.5
$ws.Cells.Item(27, 10).Value = 782.5789
$ws.Cells.Item(27, 11).Value = 417.5
$ws.Cells.Item(27, 12).Value = 782.5789
$ws.Cells.Item(27, 13).Value = -310.5
$ws.Cells.Item(27, 14).Value = -996.5789

# LTW row 123
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(123, 8).Value = 21428.572
$ws.Cells.Item(123, 10).Value = 21428.572
$ws.Cells.Item(123, 12).Value = 21428.572
$ws.Cells.Item(123, 14).Value = -31228.572

# LTW row 125
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(125, 8).Value = 20000
$ws.Cells.Item(125, 10).Value = 20000
$ws.Cells.Item(125, 12).Value = 20000
$ws.Cells.Item(125, 14).Value = -29840

# LTW row 127
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(127, 8).Value = 34275.715
$ws.Cells.Item(127, 10).Value = 34275.715
$ws.Cells.Item(127, 12).Value = 34275.715
$ws.Cells.Item(127, 14).Value = -44195.715

# LTW row 128
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(128, 8).Value = 32000
$ws.Cells.Item(128, 10).Value = 32000
$ws.Cells.Item(128, 12).Value = 32000
$ws.Cells.Item(128, 14).Value = -41960

# LTW row 130
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(130, 8).Value = 19948.092
$ws.Cells.Item(130, 10).Value = 19948.092
$ws.Cells.Item(130, 12).Value = 19948.092
$ws.Cells.Item(130, 14).Value = -29988.092

# WVR row 123
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(123, 8).Value = 22500
$ws.Cells.Item(123, 10).Value = 22500
$ws.Cells.Item(123, 12).Value = 22500
$ws.Cells.Item(123, 14).Value = -32300

# WVR row 125
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(125, 8).Value = 0
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(125, 14).Value = 0

# WVR row 128
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(128, 8).Value = 32331.889
$ws.Cells.Item(128, 10).Value = 32331.889
$ws.Cells.Item(128, 12).Value = 32331.889
$ws.Cells.Item(128, 14).Value = -42291.889
